$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shots")

# Insert a new row at row 18 (pushes the existing rows 18-26 down to 19-27,
# and Excel auto-updates the $A$26 absolute references used by column C to $A$27).
$ws.Rows.Item(18).Insert()

# New auto-shot data point (x = seconds, y = RPM).
$ws.Range("A18").Value = 7.38
$ws.Range("B18").Value = 1900

# Recalibrated RPM values for the points around the new one.
$ws.Range("B17").Value = 1835
$ws.Range("B19").Value = 1920
$ws.Range("B20").Value = 1940

# Calibrated-speed formula for the new row, matching the shared formula
# used by the rest of column C.
$ws.Range("C18").Formula = "=B18-A18/`$A`$27*350"

# Forecast helper cell added alongside the new point.
$ws.Range("D17").Formula = "=FORECAST(7.43,B18:B19,A18:A19)"
